# Spread the single BIEF39 ETF holding across two institutions:
# duplicate the existing ETF row, split the quantity 2 -> 1 + 1 between
# "INSTITUIÇÃO1" and "INSTITUIÇÃO2", and refresh the sheet's Total.

$wb = $excel.ActiveWorkbook

$wsEtf = $wb.Worksheets.Item("ETF")

# Duplicate row 2 (the BIEF39 holding) into a new row 3, carrying over
# all formatting/styles along with the values.
$wsEtf.Rows.Item(2).Copy()
$wsEtf.Rows.Item(3).Insert()

# Inserting a row pushes every row below down by one, including the
# blank formatting-only row that used to sit at the very bottom of the
# sheet (row 1048576) - real Excel has no row 1048577, so that row
# simply falls off the grid. Drop it here too so the sheet keeps a
# single trailing placeholder row instead of gaining a phantom one.
$wsEtf.Rows.Item(1048576).Delete()

# Match the original row's height on the newly inserted row.
$wsEtf.Rows.Item(3).RowHeight = $wsEtf.Rows.Item(2).RowHeight

# Row 2 now belongs to the first institution, with half the quantity.
$wsEtf.Range("B2").Value = "INSTITUIÇÃO1"
$wsEtf.Range("F2").Value = 1
$wsEtf.Range("G2").Value = 1

# Row 3 is the same product held at a second institution.
$wsEtf.Range("B3").Value = "INSTITUIÇÃO2"
$wsEtf.Range("F3").Value = 1
$wsEtf.Range("G3").Value = 1

# The Total row (previously K4/K5) shifted down to K5/K6; refresh the
# summed "Valor Atualizado" total to reflect both institutions' rows.
$wsEtf.Range("K6").Value = 4.44

# Refresh the selection on the other sheets (matches where the user's
# cursor ended up while reviewing the change) and finally land on / make
# active the ETF sheet.
$wsAcoes = $wb.Worksheets.Item("Ações")
[void]$wsAcoes.Activate()
[void]$wsAcoes.Range("A3").Select()

$wsFundo = $wb.Worksheets.Item("Fundo de Investimento")
[void]$wsFundo.Activate()
[void]$wsFundo.Range("L5").Select()

[void]$wsEtf.Activate()
[void]$wsEtf.Range("A3").Select()
